$wb = $excel.ActiveWorkbook

# Sheet1, Sheet2, Sheet3 each have a "Current Date:" label in A3 and a date
# value in A4. The demo no longer stamps the live system date/time
# (sy-datum/sy-uzeit), so the label becomes just "Date:" and the value
# becomes a small static number instead of a captured date serial.
foreach ($name in @("Sheet1", "Sheet2", "Sheet3")) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("A3").Value = "Date:"
    $ws.Range("A4").Value = 57
}
